$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price and 1h volume %) per upstream diff.
$changes = [ordered]@{
    "D2" = '281.48'
    "E2" = '1.54%'
    "D3" = '27.63'
    "D4" = '4.842'
    "E4" = '-1.18%'
    "D5" = '0.06513'
    "E5" = '1.42%'
    "D6" = '7.104'
    "E6" = '2.24%'
    "D7" = '1.272'
    "E7" = '2.72%'
    "D8" = '0.9199'
    "E8" = '4.24%'
    "D9" = '0.1556'
    "E9" = '2.33%'
    "D10" = '0.06444'
    "E10" = '28.43%'
    "D11" = '0.07513'
    "E11" = '0.11%'
    "E12" = '0.36%'
    "D13" = '0.09020'
    "E13" = '0.15%'
    "D14" = '0.001576'
    "E14" = '0.61%'
    "D15" = '0.0006407'
    "E15" = '-0.30%'
    "D16" = '0.005999'
    "E16" = '5.09%'
    "D17" = '3.494'
    "E17" = '0.98%'
    "D18" = '3.316'
    "E18" = '0.04%'
    "D19" = '2.225'
    "E19" = '-2.07%'
    "D20" = '0.3144'
    "E20" = '0.26%'
    "D21" = '0.1342'
    "E21" = '0.38%'
    "D22" = '3.887'
    "E22" = '-0.54%'
    "D23" = '0.04407'
    "E23" = '-0.78%'
    "D24" = '0.1487'
    "E24" = '7.71%'
    "D25" = '0.001167'
    "E25" = '-0.64%'
    "D26" = '0.004394'
    "E26" = '13.80%'
    "D28" = '0.0001175'
    "E28" = '-2.17%'
    "D29" = '0.0001640'
    "E29" = '-15.29%'
    "D40" = '0.04110'
    "E40" = '-0.83%'
    "D41" = '0.1397'
    "E41" = '18.75%'
    "D42" = '0.006301'
    "E42" = '-7.54%'
    "D43" = '0.002082'
    "E43" = '-12.99%'
    "D44" = '0.01162'
    "E44" = '-0.74%'
    "D45" = '0.00005487'
    "E45" = '5.45%'
    "D46" = '1.628'
    "E46" = '9.53%'
    "D47" = '0.01830'
    "E47" = '-9.63%'
}

foreach ($cellRef in $changes.Keys) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    # Leading apostrophe forces text entry so numeric-looking / percent-looking
    # strings are stored verbatim as text (matching the original inlineStr cells)
    # instead of being coerced into Excel numbers.
    $rng.Value = "'" + $changes[$cellRef]
    # Restore the original (default) cell style, since the quote-prefix entry
    # above would otherwise tag the cell with a new "quotePrefix" style record.
    $rng.Style = $origStyle
}

Write-Output ("Updated {0} cells" -f $changes.Count)
